# Applies the "Add files via upload" edit to the sales-revenue workbook:
#   - Customer_Reviews (column C, rows 2-51) values are rounded from their
#     long-decimal raw form to whole numbers.
#   - Columns A:D are best-fit/auto-fitted to their (now shorter) content.
#   - The thin box border that outlined the bold header row (A1:D1) is removed.
#   - The saved selection moves to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round the Customer_Reviews column (C) to the nearest whole number ---
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = [Math]::Round($cell.Value(), 0)
}

# --- Best-fit the column widths now that the header/data footprint changed ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# --- Strip the thin outline border that was around the header cells ---
$ws.Range("A1:D1").Borders.LineStyle = -4142

# --- Leave the selection where the workbook was saved (D5) ---
$ws.Range("D5").Select()
